# Add a new "asesoria" (advisory) record for student "david soto" to the
# "asesorias" sheet, mirroring the existing row for "Sebastian Palacio" /
# "Juan Carlos Gil" (same advisor, motive and time slot) but with a new
# student name and date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asesorias")
$ws.Activate()

$ws.Range("A5").Value = "david soto"
$ws.Range("B5").Value = "s"
$ws.Range("C5").Value = "Juan Carlos Gil"
$ws.Range("D5").Value = "Consulta general"
$ws.Range("E5").Value = "24-11-2023"
$ws.Range("F5").Value = "00:20 - 00:40"

# Reflect the newly-entered row as the active selection, matching the
# author's on-screen state after typing the new record.
$ws.Range("A5:F5").Select()
